# EditAPI-style update: rename circle tabs and refresh the report date,
# then leave the workbook open on the "Circle 11" tab (matches last edit).

$wb = $excel.ActiveWorkbook

# Drop the "(n)" duplicate-tab suffixes that Excel appended when these
# sheets were copied from the Circle 7 template.
$wb.Worksheets.Item("Circle 7 (2)").Name = "Circle 7"
$wb.Worksheets.Item("Circle 8 (3)").Name = "Circle 8"
$wb.Worksheets.Item("Circle 9 (4)").Name = "Circle 9"
$wb.Worksheets.Item("Circle 10 (5)").Name = "Circle 10"

# Update the report date on the Circle 6 sheet.
$ws6 = $wb.Worksheets.Item("Circle 6")
$ws6.Range("K2").Value = "Date: 11-09-2020"

# Leave the workbook focused on the Circle 11 tab.
$wb.Worksheets.Item("Circle 11").Activate()
